$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster had Ja Morant in row 2 and Tyler Herro in row 14. Swap the two
# rows' contents (player, position, team) so Tyler Herro ends up in row 2
# and Ja Morant ends up in row 14; every other row is untouched.

$a2  = $ws.Range("A2").Text
$b2  = $ws.Range("B2").Text
$c2  = $ws.Range("C2").Text
$a14 = $ws.Range("A14").Text
$b14 = $ws.Range("B14").Text
$c14 = $ws.Range("C14").Text

$ws.Range("A2").Value = $a14
$ws.Range("B2").Value = $b14
$ws.Range("C2").Value = $c14

$ws.Range("A14").Value = $a2
$ws.Range("B14").Value = $b2
$ws.Range("C14").Value = $c2
